$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the existing date cell (A2) onto the new
# date cell A4 first, so no new number-format/style entries get created.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row of data (row 4): date, duration, description
$ws.Range("A4").Value = (Get-Date -Year 2016 -Month 12 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B4").Value = "4h"
$ws.Range("C4").Value = "Einarbeitung/Tutorials anschauen bezüglich asp.net core"

$ws.Range("A4").Select()
